# ===========================================================================
# WordPress Plugin Proposal -- apply the "WP plugin proposal complete" edit.
#
# Strategy: each touched paragraph's Range is rewritten in one shot via
# Range.InsertXML with the exact target run/proofErr structure (this lets us
# control run boundaries and <w:proofErr> spell-check markers precisely,
# which plain Find/Replace can't do). New trailing paragraphs are created
# with Range.InsertParagraphAfter() and then populated the same way.
# ===========================================================================

$d = $word.ActiveDocument
$emptyPXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'

function Assert-ParaText($index, $expected) {
    $actual = $d.Paragraphs.Item($index).Range.Text.TrimEnd([char]13, [char]7)
    if ($actual -ne $expected) {
        throw "Paragraph $index text mismatch. Expected [$expected] got [$actual]"
    }
}

# ---------------------------------------------------------------------------
# 1. Title paragraph: "WordPress Plugin Proposal"
#    -> split into 3 runs, "Plugin" wrapped in proofErr spellStart/spellEnd
# ---------------------------------------------------------------------------
Assert-ParaText 1 "WordPress Plugin Proposal"
$titleXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="40"/></w:rPr><w:t xml:space="preserve">WordPress </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="40"/></w:rPr><w:t>Plugin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="40"/></w:rPr><w:t xml:space="preserve"> Proposal</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(1).Range.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 2. Intro paragraph (was "I would like to create a plugin for WordPress...")
#    -> fully rewritten with many runs / proofErr wraps
# ---------------------------------------------------------------------------
Assert-ParaText 3 "I would like to create a plugin for WordPress to link social media pages. This plugin will give direct access to the selected social media page along with a “Follow” button for the user to easily start following the page directly from the WordPress site that the plugin is used on. This will prevent the user from having to leave the WordPress page."
$introXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>This proposal is for the creation of</w:t></w:r><w:r><w:t xml:space="preserve"> a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>plugin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for WordP</w:t></w:r><w:r><w:t>ress to link social media pages to sidebars on a WordPress page</w:t></w:r><w:r><w:t xml:space="preserve"> This </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>plugin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> will give</w:t></w:r><w:r><w:t xml:space="preserve"> users a link to</w:t></w:r><w:r><w:t xml:space="preserve"> direct access of</w:t></w:r><w:r><w:t xml:space="preserve"> the selected social media page along with a &#8220;Follow&#8221; button for the user to easily start following the page directly from the WordPress site that the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>plugin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is used on. This will prevent the user from having to leave the WordPress page.</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(3).Range.InsertXML($introXml)

# ---------------------------------------------------------------------------
# 3. Timeline-intro paragraph: "This " + "project should take 2" stay as-is;
#    the trailing run gets rewritten/split with proofErr wraps.
# ---------------------------------------------------------------------------
Assert-ParaText 5 "This project should take 2 weeks for coding, testing, and debugging. Once completed, the admin should have the ability to enter the url of their Facebook, Twitter, Instagram, Tumblr, YouTube, Vimeo, Google +, and/or Pintrest page. A subscribe or follow button will be available to add to that."
$timelineXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">This </w:t></w:r><w:r><w:t>project should take 2</w:t></w:r><w:r><w:t xml:space="preserve"> weeks for coding, testing, and debugging. Once completed, the admin should have the ability to enter the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> of their </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Facebook</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, Twitter, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Instagram</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tumblr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, YouTube, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Vimeo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, Google +, and/or </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Pintrest</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> page. A subscribe or follow button will be available to add to that.</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(5).Range.InsertXML($timelineXml)

# ---------------------------------------------------------------------------
# 4. Append new paragraphs after "Nov 16-22 - test and debug" (paragraph 13):
#    3 blank paragraphs, a COMMENTS banner, a blank paragraph, and a
#    3-run comments paragraph.
# ---------------------------------------------------------------------------
Assert-ParaText 13 "Nov 16-22 – test and debug"

[void]$d.Paragraphs.Item(13).Range.InsertParagraphAfter()    # -> paragraph 14 (blank)
[void]$d.Paragraphs.Item(14).Range.InsertXML($emptyPXml)

[void]$d.Paragraphs.Item(14).Range.InsertParagraphAfter()    # -> paragraph 15 (blank)
[void]$d.Paragraphs.Item(15).Range.InsertXML($emptyPXml)

[void]$d.Paragraphs.Item(15).Range.InsertParagraphAfter()    # -> paragraph 16 (blank)
[void]$d.Paragraphs.Item(16).Range.InsertXML($emptyPXml)

[void]$d.Paragraphs.Item(16).Range.InsertParagraphAfter()    # -> paragraph 17 (COMMENTS banner)
$commentsHeaderXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>************COMMENTS*************</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(17).Range.InsertXML($commentsHeaderXml)

[void]$d.Paragraphs.Item(17).Range.InsertParagraphAfter()    # -> paragraph 18 (blank)
[void]$d.Paragraphs.Item(18).Range.InsertXML($emptyPXml)

[void]$d.Paragraphs.Item(18).Range.InsertParagraphAfter()    # -> paragraph 19 (final comments paragraph)
$finalCommentsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>I&#8217;m interested to get started on this one and see how plugin</w:t></w:r><w:r><w:t>s work. I&#8217;m not sure if this idea</w:t></w:r><w:r><w:t xml:space="preserve"> is too ambitious but I think it would be a good one to start with.</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(19).Range.InsertXML($finalCommentsXml)

Write-Output ("Done. Final paragraph count: " + $d.Paragraphs.Count)
